# Auto-generated Excel COM-interop edit script
# Updates crypto price (D) and 1h volume change (E) columns, plus
# coin name/link for row 50 (RenderToken -> Maker), per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) contain formatted numeric-looking text
# (e.g. '59.450.03', '1.00', '0.340') that must stay literal text,
# so force Text number format before assigning the value.
$priceCells = @('D2', 'D3', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D12', 'D14', 'D15', 'D16', 'D18', 'D19', 'D20', 'D21', 'D22', 'D24', 'D25', 'D26', 'D27', 'D29', 'D31', 'D32', 'D33', 'D35', 'D36', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50')
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = '@'
}

# Price (D) updates
$ws.Range('D2').Value = '59.450.03'
$ws.Range('D3').Value = '2.640.16'
$ws.Range('D5').Value = '516.62'
$ws.Range('D6').Value = '147.17'
$ws.Range('D7').Value = '0.998'
$ws.Range('D8').Value = '0.572'
$ws.Range('D9').Value = '2.670.75'
$ws.Range('D10').Value = '6.47'
$ws.Range('D12').Value = '0.340'
$ws.Range('D14').Value = '3.137.50'
$ws.Range('D15').Value = '59.461.12'
$ws.Range('D16').Value = '21.26'
$ws.Range('D18').Value = '2.671.16'
$ws.Range('D19').Value = '4.61'
$ws.Range('D20').Value = '345.99'
$ws.Range('D21').Value = '10.53'
$ws.Range('D22').Value = '6.21'
$ws.Range('D24').Value = '61.47'
$ws.Range('D25').Value = '0.424'
$ws.Range('D26').Value = '2.769.58'
$ws.Range('D27').Value = '0.995'
$ws.Range('D29').Value = '0.0₃0822'
$ws.Range('D31').Value = '1.00'
$ws.Range('D32').Value = '6.52'
$ws.Range('D33').Value = '19.11'
$ws.Range('D35').Value = '149.99'
$ws.Range('D36').Value = '1.05'
$ws.Range('D38').Value = '1.17'
$ws.Range('D39').Value = '0.872'
$ws.Range('D40').Value = '36.68'
$ws.Range('D41').Value = '3.73'
$ws.Range('D42').Value = '1.43'
$ws.Range('D43').Value = '286.42'
$ws.Range('D45').Value = '0.0995'
$ws.Range('D46').Value = '0.995'
$ws.Range('D47').Value = '19.81'
$ws.Range('D48').Value = '0.0545'
$ws.Range('D49').Value = '0.0233'
$ws.Range('D50').Value = '1.995.63'

# Coin / Link / Volume(1h) updates
$ws.Range('E2').Value = '  -1.18%  '
$ws.Range('E3').Value = '  +0.79%  '
$ws.Range('E4').Value = '  +0.34%  '
$ws.Range('E5').Value = '  -0.96%  '
$ws.Range('E6').Value = '  -1.08%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  +0.37%  '
$ws.Range('E9').Value = '  +1.81%  '
$ws.Range('E10').Value = '  +2.61%  '
$ws.Range('E11').Value = '  +1.25%  '
$ws.Range('E12').Value = '  -0.46%  '
$ws.Range('E13').Value = '  -1.47%  '
$ws.Range('E14').Value = '  +2.04%  '
$ws.Range('E15').Value = '  -1.15%  '
$ws.Range('E16').Value = '  +0.32%  '
$ws.Range('E17').Value = '  +0.39%  '
$ws.Range('E18').Value = '  +1.82%  '
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('E20').Value = '  +1.31%  '
$ws.Range('E21').Value = '  +1.03%  '
$ws.Range('E22').Value = '  +1.59%  '
$ws.Range('E23').Value = '  +0.58%  '
$ws.Range('E24').Value = '  +1.50%  '
$ws.Range('E25').Value = '  +1.01%  '
$ws.Range('E26').Value = '  +1.34%  '
$ws.Range('E27').Value = '  -0.48%  '
$ws.Range('E28').Value = '  -0.08%  '
$ws.Range('E29').Value = '  +1.54%  '
$ws.Range('E30').Value = '  +2.10%  '
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('E32').Value = '  +8.71%  '
$ws.Range('E33').Value = '  +0.80%  '
$ws.Range('E34').Value = '  -0.10%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('E36').Value = '  +14.38%  '
$ws.Range('E37').Value = '  +2.78%  '
$ws.Range('E38').Value = '  +2.99%  '
$ws.Range('E39').Value = '  +1.02%  '
$ws.Range('E40').Value = '  +0.63%  '
$ws.Range('E41').Value = '  +3.13%  '
$ws.Range('E42').Value = '  -0.36%  '
$ws.Range('E43').Value = '  -1.00%  '
$ws.Range('E44').Value = '  -0.85%  '
$ws.Range('E45').Value = '  -0.96%  '
$ws.Range('E46').Value = '  -0.32%  '
$ws.Range('E47').Value = '  +1.84%  '
$ws.Range('E48').Value = '  -0.38%  '
$ws.Range('E49').Value = '  +0.74%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('E50').Value = '  +1.62%  '
$ws.Range('E51').Value = '  -1.23%  '
